$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as genuine TEXT (shared-string) content, even when
# the text looks like a number (e.g. "-0", "-3.6"), without Excel's normal
# numeric auto-detection converting it into a numeric cell and without
# leaving behind a stray cell style (NumberFormat="@" would add one).
# We do this the way Excel itself would if you built the value with a
# formula and then froze it: enter ="<text>"  then Copy + Paste-Values.
function Set-TextValue($cell, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)   # xlPasteValues
}

# New column G ("dSoH (ppm)") is appended after the existing data -- F keeps
# meaning "Adj. rev." exactly as before, only its row3/row4 numbers change.
$ws.Cells.Item(2,7).Value = "dSoH (ppm)"

# Row-1 index label for the new column.
Set-TextValue $ws.Cells.Item(1,6) "5"

# Updated "Adj. rev." figures (col F) and their source "Trading rev." (col C)
# / "HW cost" (col D) figures.
Set-TextValue $ws.Cells.Item(3,3) "-4.1"
Set-TextValue $ws.Cells.Item(3,6) "-4.1"
Set-TextValue $ws.Cells.Item(4,3) "-3"
Set-TextValue $ws.Cells.Item(4,4) "-0.6"
Set-TextValue $ws.Cells.Item(4,6) "-3.6"

# New "dSoH (ppm)" values.
Set-TextValue $ws.Cells.Item(3,7) "56.1"
Set-TextValue $ws.Cells.Item(4,7) "50.3"

# Give the new column a sensible best-fit width, matching its neighbours
# (engine quantises to whole pixels, so 8.5 is the closest match to the
# ~8.457 authored width).
$ws.Columns(6).ColumnWidth = 7.6
